$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function Get-ParagraphIndexContaining([string]$needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

function Set-ParagraphXmlByNeedle([string]$needle, [string]$innerXml) {
    $idx = Get-ParagraphIndexContaining $needle
    if ($idx -lt 0) {
        throw "Paragraph containing '$needle' not found"
    }
    # Re-fetch a fresh Range right before InsertXML - a Range that has been
    # narrowed via Find.Execute leaves stray zero-width markers (w:proofErr)
    # behind on replace, whereas a pristine Paragraph.Range swaps the whole
    # paragraph content cleanly.
    $range = $d.Paragraphs.Item($idx).Range
    $xml = '<w:p xmlns:w="' + $wNs + '">' + $innerXml + '</w:p>'
    $range.InsertXML($xml)
}

# --- "What is needed:" -> merge the two runs, drop proofErr gramStart/gramEnd ---
Set-ParagraphXmlByNeedle "What is " '<w:r><w:t>What is needed:</w:t></w:r>'

# --- "Lines with antialiasing (possibly using shaders)" -> merge runs, drop proofErr ---
Set-ParagraphXmlByNeedle "Lines with " (
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Lines with antialiasing (possibly using shaders)</w:t></w:r>'
)

# --- "{ Postprocessing, fragment shaders (for effects) }" -> merge middle runs, drop proofErr ---
Set-ParagraphXmlByNeedle "Postprocessing" (
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
    '<w:rPr><w:color w:val="7F7F7F" w:themeColor="text1" w:themeTint="80"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:color w:val="7F7F7F" w:themeColor="text1" w:themeTint="80"/></w:rPr><w:t xml:space="preserve">{ </w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="7F7F7F" w:themeColor="text1" w:themeTint="80"/></w:rPr><w:t>Postprocessing, fragment shaders (for effects)</w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="7F7F7F" w:themeColor="text1" w:themeTint="80"/></w:rPr><w:t xml:space="preserve"> }</w:t></w:r>'
)

# --- "Support for iPhone 4* and iPad* resolutions (960x640, 1024x768)" -> merge runs, drop proofErr ---
Set-ParagraphXmlByNeedle "Support for iPhone" (
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Support for iPhone 4* and iPad*</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> resolutions (960x640, 1024x768)</w:t></w:r>'
)

# --- "Time limit increase { time_limit [seconds] }" -> merge runs, drop proofErr ---
Set-ParagraphXmlByNeedle "Time limit increase" (
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Time limit increase</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> { time_limit [seconds] }</w:t></w:r>'
)

# --- "Bomb, bonus items increase { bonus_freq }" -> merge runs, drop proofErr ---
Set-ParagraphXmlByNeedle "Bomb, bonus items increase" (
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Bomb, bonus items increase</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> { bonus_freq }</w:t></w:r>'
)

# --- "Target charges increase { target_charges }" -> merge runs, drop proofErr ---
Set-ParagraphXmlByNeedle "Target charges increase" (
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Target charges increase</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> { target_charges }</w:t></w:r>'
)

# --- "Worm starts to appear { worm_level }" -> merge runs, drop proofErr ---
Set-ParagraphXmlByNeedle "Worm starts to appear" (
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Worm starts to appear</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> { worm_level }</w:t></w:r>'
)

# --- "Worm increases length { worm_length }" -> merge runs, drop proofErr ---
Set-ParagraphXmlByNeedle "Worm increases length" (
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Worm increases length</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> { worm_length }</w:t></w:r>'
)

# --- "Worm starts to rotate tiles with increasing frequency { worm_tile_rotations }" -> merge runs, drop proofErr ---
Set-ParagraphXmlByNeedle "Worm starts to rotate tiles" (
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Worm starts to rotate tiles with increasing frequency</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> { worm_tile_rotations }</w:t></w:r>'
)

# --- "Difficulty (affects level progression) { ? }" -> merge "{ ?" run with " }" run, drop proofErr ---
Set-ParagraphXmlByNeedle "Difficulty" (
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Difficulty</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> (affects level progression)</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="7F7F7F" w:themeColor="text1" w:themeTint="80"/></w:rPr><w:t>{ ? }</w:t></w:r>'
)

# --- "Move worm management to GameTable class" -> drop proofErr around GameTable (runs unchanged) ---
Set-ParagraphXmlByNeedle "Move worm management" (
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Move worm management to </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/></w:rPr><w:t>GameTable</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> class</w:t></w:r>'
)

# --- "Bonus items { CHARGES, bombs, clocks } put them into GameTable class" -> drop proofErr around GameTable ---
Set-ParagraphXmlByNeedle "Bonus items { CHARGES" (
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Bonus items</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> { CHARGES, bombs, clocks }</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> put them into </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/></w:rPr><w:t>GameTable</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> class</w:t></w:r>'
)

# --- Add a new bullet right after "Textured polygons for lightning" ---
$idx = Get-ParagraphIndexContaining "Textured polygons for lightning"
if ($idx -lt 0) {
    throw "Paragraph 'Textured polygons for lightning' not found"
}
$endOfPara = $d.Paragraphs.Item($idx).Range.End
$insertAt = $d.Range($endOfPara, $endOfPara)
$newParaXml = '<w:p xmlns:w="' + $wNs + '">' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Splash screen with tiles { PIPE CHARGER }</w:t></w:r>' +
    '</w:p>'
$insertAt.InsertXML($newParaXml)

Write-Output "DONE"
